$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "A80" "7304"
Set-TextValue "B80" "9/24/2025"
$ws.Range("C80").Value = "VALLE 796"
Set-TextValue "D80" "6"
Set-TextValue "E80" "809979725"
$ws.Range("F80").Value = "PEBCOM"
$ws.Range("G80").Value = "Pendiente"
$ws.Range("H80").Value = "Cambiar columna colocar la nueva fuera del cantero"
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = "Cambio"
$ws.Range("K80").Value = "Sin equipos"
$ws.Range("L80").Value = "Pasante"
$ws.Range("M80").Value = -58.439096
$ws.Range("N80").Value = -34.624889
$ws.Range("O80").Value = "Boedo"
$ws.Range("P80").Value = "Capital Sur"
$ws.Range("Q80").Value = "PCH-C"
$ws.Range("R80").Value = "Fuera de Poligono OVL"

Set-TextValue "A81" "7317"
Set-TextValue "B81" "9/25/2025"
$ws.Range("C81").Value = "MARMOL, JOSE 588"
Set-TextValue "D81" "5"
Set-TextValue "E81" "809979740"
$ws.Range("F81").Value = "PEBCOM"
$ws.Range("G81").Value = "Pendiente"
$ws.Range("H81").Value = "Picada"
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = "Cambio"
$ws.Range("K81").Value = "Sin equipos"
$ws.Range("L81").Value = "Pasante"
$ws.Range("M81").Value = -58.425357
$ws.Range("N81").Value = -34.620223
$ws.Range("O81").Value = "Almagro"
$ws.Range("P81").Value = "Capital Sur"
$ws.Range("Q81").Value = "ALM-B"
$ws.Range("R81").Value = "Fuera de Poligono OVL"
